# Remove the duplicate "현대힐스테이트" tracking entry (A=13698, 84) that
# previously lived at row 229 - it is an exact duplicate of row 216.
# Deleting the whole row shifts every subsequent row up by one, which is
# exactly what the target diff shows (old rows 230-241 become new rows
# 229-240, and the sheet's used range shrinks from C241 to C240).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("229:229").Delete()

# Update the view state to match where the author had scrolled/selected
# after making the edit.
$ws.Range("K227:L227").Select()
$excel.ActiveWindow.ScrollRow = 208
$excel.ActiveWindow.ScrollColumn = 1
